# Update the NATMI LR-pairs worksheet (Ntf5-Ngfr) with the re-run (new TPM) results.
# The sending/target cluster set grew from {FAPs, Resolving-Mac} to
# {FAPs, MuSCs, Resolving-Mac}, so the 8-row table becomes a 9-row table
# (3 senders x 3 targets) and all computed metrics are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the extra data row (old data occupied rows 2-9; new data needs rows 2-10).
$ws.Rows("9:9").Insert()

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @("FAPs", "Ntf5", "Ngfr", "ECs", 1, 0.3333333333333333, 0.563831, 1.691493, 0.5426742997499485, 0.5426742997499486, 1, 0.3333333333333333, 0.261293, 0.783879, 0.0361900776512412, 0.03619007765124121, 0.147325093483, 1.325925841347, 0.01963942504728358, 0.01963942504728359),
    @("FAPs", "Ntf5", "Ngfr", "FAPs", 1, 0.3333333333333333, 0.563831, 1.691493, 0.5426742997499485, 0.5426742997499486, 3, 1, 0.7775033333333333, 2.33251, 0.1076871787894517, 0.1076871787894517, 0.4383804819366666, 3.94542433743, 0.05843906434161322, 0.05843906434161324),
    @("FAPs", "Ntf5", "Ngfr", "MuSCs", 1, 0.3333333333333333, 0.563831, 1.691493, 0.5426742997499485, 0.5426742997499486, 3, 1, 6.181221333333333, 18.543664, 0.8561227435593071, 0.8561227435593072, 3.485164205594666, 31.366477850352, 0.4645958103610517, 0.4645958103610519),
    @("MuSCs", "Ntf5", "Ngfr", "ECs", 1, 0.3333333333333333, 0.05363299999999999, 0.160899, 0.05162052231695133, 0.05162052231695134, 1, 0.3333333333333333, 0.261293, 0.783879, 0.0361900776512412, 0.03619007765124121, 0.014013927469, 0.126125347221, 0.001868150711048098, 0.001868150711048099),
    @("MuSCs", "Ntf5", "Ngfr", "FAPs", 1, 0.3333333333333333, 0.05363299999999999, 0.160899, 0.05162052231695133, 0.05162052231695134, 3, 1, 0.7775033333333333, 2.33251, 0.1076871787894517, 0.1076871787894517, 0.04169983627666666, 0.37529852649, 0.00555886841595042, 0.005558868415950422),
    @("MuSCs", "Ntf5", "Ngfr", "MuSCs", 1, 0.3333333333333333, 0.05363299999999999, 0.160899, 0.05162052231695133, 0.05162052231695134, 3, 1, 6.181221333333333, 18.543664, 0.8561227435593071, 0.8561227435593072, 0.3315174437706666, 2.983656993936, 0.04419350318995281, 0.04419350318995283),
    @("Resolving-Mac", "Ntf5", "Ngfr", "ECs", 2, 0.6666666666666666, 0.421522, 1.264566, 0.4057051779331001, 0.4057051779331002, 1, 0.3333333333333333, 0.261293, 0.783879, 0.0361900776512412, 0.03619007765124121, 0.110140747946, 0.9912667315139999, 0.01468250189290952, 0.01468250189290953),
    @("Resolving-Mac", "Ntf5", "Ngfr", "FAPs", 2, 0.6666666666666666, 0.421522, 1.264566, 0.4057051779331001, 0.4057051779331002, 3, 1, 0.7775033333333333, 2.33251, 0.1076871787894517, 0.1076871787894517, 0.3277347600733333, 2.94961284066, 0.04368924603188807, 0.04368924603188809),
    @("Resolving-Mac", "Ntf5", "Ngfr", "MuSCs", 2, 0.6666666666666666, 0.421522, 1.264566, 0.4057051779331001, 0.4057051779331002, 3, 1, 6.181221333333333, 18.543664, 0.8561227435593071, 0.8561227435593072, 2.605520778869333, 23.449687009824, 0.3473334300083025, 0.3473334300083026)
)

$rowIndex = 2
foreach ($rowValues in $data) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = $columns[$i] + $rowIndex
        $ws.Range($cellRef).Value = $rowValues[$i]
    }
    $rowIndex++
}
